$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "70.382.67", "0.998").
# Excel auto-converts such strings to real numbers when assigned via
# .Value, which would corrupt the text (rounding, scientific notation,
# loss of trailing zeros, loss of the "thousands dot" grouping, etc).
# Force the column to Text format first so values are stored verbatim,
# then restore the original (default) style so no visible formatting
# change is introduced.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "70.382.67"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "3.570.92"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "611.62"
$ws.Range("E5").Value = "  +4.11%  "
$ws.Range("D6").Value = "187.95"
$ws.Range("E6").Value = "  +2.20%  "
$ws.Range("D7").Value = "3.565.92"
$ws.Range("E7").Value = "  +1.70%  "
$ws.Range("D8").Value = "0.620"
$ws.Range("E8").Value = "  +1.44%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "0.214"
$ws.Range("E10").Value = "  +8.65%  "
$ws.Range("D11").Value = "0.649"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D12").Value = "54.18"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").Value = "0.0000310"
$ws.Range("E13").Value = "  +2.21%  "
$ws.Range("D14").Value = "9.48"
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").Value = "4.133.83"
$ws.Range("E15").Value = "  +1.44%  "
$ws.Range("D16").Value = "70.425.83"
$ws.Range("E16").Value = "  +1.04%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value = "12.77"
$ws.Range("E17").Value = "  +3.51%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.561.77"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").Value = "19.01"
$ws.Range("E19").Value = "  -1.49%  "
$ws.Range("D20").Value = "571.94"
$ws.Range("E20").Value = "  +7.50%  "
$ws.Range("D21").Value = "0.121"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").Value = "17.54"
$ws.Range("E23").Value = "  -2.45%  "
$ws.Range("D24").Value = "4.79"
$ws.Range("E24").Value = "  +4.42%  "
$ws.Range("D25").Value = "4.91"
$ws.Range("E25").Value = "  +1.77%  "
$ws.Range("D26").Value = "94.04"
$ws.Range("E26").Value = "  -1.59%  "
$ws.Range("D27").Value = "2.95"
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").Value = "10.98"
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("D29").Value = "9.42"
$ws.Range("E29").Value = "  +3.90%  "
$ws.Range("D30").Value = "32.46"
$ws.Range("E30").Value = "  +1.07%  "
$ws.Range("D31").Value = "7.11"
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("D32").Value = "12.25"
$ws.Range("E32").Value = "  -1.18%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.115"
$ws.Range("E33").Value = "  +2.15%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "64.41"
$ws.Range("E34").Value = "  +0.38%  "
$ws.Range("D35").Value = "3.78"
$ws.Range("E35").Value = "  +21.06%  "
$ws.Range("E36").Value = "  +3.62%  "
$ws.Range("B37").Value = "TheGraph"
$ws.Range("C37").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D37").Value = "0.406"
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "3.691.22"
$ws.Range("E38").Value = "  +9.65%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "523.15"
$ws.Range("E39").Value = "  -4.27%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "37.82"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").Value = "0.0₃0787"
$ws.Range("E42").Value = "  +3.72%  "
$ws.Range("D43").Value = "3.53"
$ws.Range("E43").Value = "  +3.66%  "
$ws.Range("E44").Value = "  +2.76%  "
$ws.Range("D45").Value = "0.0457"
$ws.Range("E45").Value = "  +4.11%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "3.49"
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Value = "2.96"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("D48").Value = "0.139"
$ws.Range("E48").Value = "  +3.12%  "
$ws.Range("E49").Value = "  +2.82%  "
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "136.67"
$ws.Range("E51").Value = "  -0.04%  "

$dRange.Style = "Normal"
